$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.119.54'
$ws.Range('E2').Value = '  -3.45%  '
$ws.Range('D3').Value = '3.140.62'
$ws.Range('E3').Value = '  -3.39%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '610.81'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.87'
$ws.Range('E6').Value = '  -6.61%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '3.140.16'
$ws.Range('E8').Value = '  -3.42%  '
$ws.Range('E9').Value = '  -3.66%  '
$ws.Range('E10').Value = '  -6.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.35'
$ws.Range('E11').Value = '  -7.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.472'
$ws.Range('E12').Value = '  -5.24%  '
$ws.Range('E13').Value = '  -7.25%  '
$ws.Range('E14').Value = '  -9.29%  '
$ws.Range('D15').Value = '3.659.44'
$ws.Range('E15').Value = '  -3.21%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.115'
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '64.117.84'
$ws.Range('E17').Value = '  -3.66%  '
$ws.Range('D18').Value = '3.143.40'
$ws.Range('E18').Value = '  -3.16%  '
$ws.Range('E19').Value = '  -8.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '476.73'
$ws.Range('E20').Value = '  -5.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.75'
$ws.Range('E21').Value = '  -4.51%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.704'
$ws.Range('E22').Value = '  -6.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.77'
$ws.Range('E23').Value = '  -4.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.62'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.32'
$ws.Range('E25').Value = '  -4.42%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.85'
$ws.Range('E27').Value = '  -5.70%  '
$ws.Range('E28').Value = '  -7.84%  '
$ws.Range('E29').Value = '  -8.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.73'
$ws.Range('E30').Value = '  -3.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.113'
$ws.Range('E31').Value = '  -10.10%  '
$ws.Range('E32').Value = '  -5.32%  '
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.17'
$ws.Range('E34').Value = '  -6.35%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.12'
$ws.Range('E35').Value = '  -2.46%  '
$ws.Range('E36').Value = '  -7.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.42'
$ws.Range('E37').Value = '  -3.71%  '
$ws.Range('E38').Value = '  -5.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '462.34'
$ws.Range('E39').Value = '  -6.67%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.91'
$ws.Range('E40').Value = '  -12.54%  '
$ws.Range('E41').Value = '  -6.50%  '
$ws.Range('E42').Value = '  -8.15%  '
$ws.Range('E43').Value = '  -4.79%  '
$ws.Range('D44').Value = '2.847.66'
$ws.Range('E44').Value = '  -4.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.266'
$ws.Range('E45').Value = '  -9.21%  '
$ws.Range('E46').Value = '  -10.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.41'
$ws.Range('E47').Value = '  -8.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.999'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.37'
$ws.Range('E49').Value = '  -6.09%  '
$ws.Range('E50').Value = '  -4.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '118.35'
$ws.Range('E51').Value = '  -2.13%  '
